$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Usuario" (username) test value in A2
$ws.Range("A2").Value = "testeemnswd"

# Update the "Senha" (password) test value in C2. The cell's formatting was
# also reset to the default style (its font no longer uses the old custom
# font) as part of this edit.
$ws.Range("C2").Style = "Normal"
$ws.Range("C2").Value = "Teste@2022"

# Restore the active selection to C2, matching the saved sheet view state.
$ws.Range("C2").Select()
